$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-11) need to be reordered (e.g. re-sorted by date/week).
# Capture the full original row contents (columns A:R) before overwriting
# anything, then write them back out in the new order.

$orig = @{}
for ($r = 2; $r -le 11; $r++) {
    $orig[$r] = $ws.Range("A$r`:R$r").Value2
}

# Mapping: new row number -> original row number it should now contain
$order = @{
    2  = 11
    3  = 9
    4  = 7
    5  = 4
    6  = 8
    7  = 5
    8  = 6
    9  = 2
    10 = 3
    11 = 10
}

foreach ($newRow in 2..11) {
    $oldRow = $order[$newRow]
    $ws.Range("A$newRow`:R$newRow").Value = $orig[$oldRow]
}
